{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n\nconst body = context.document.body;\n\n// 1) Delete the first paragraph, which contains only the inline picture\n//    (Picture 1 / 34.jpg). It has no text, so find it positionally.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length > 0) {\n  const firstPara = paragraphs.items[0];\n  firstPara.load(\"text\");\n  const pics = firstPara.inlinePictures;\n  pics.load(\"items\");\n  await context.sync();\n\n  if (pics.items.length > 0) {\n    firstPara.delete();\n    await context.sync();\n  }\n}\n\n// 2) Fix the typo \"contalner\" -> \"container\" in \"... was slowly heated up.\"\nconst typoResults = body.search(\"contalner was slowly heated up.\", { matchCase: true });\ntypoResults.load(\"items\");\nawait context.sync();\ntypoResults.items.forEach((r) => {\n  r.insertText(\"container was slowly heated up.\", Word.InsertLocation.replace);\n});\nawait context.sync();\n\n// 3) \"PA\" -> \"fA BT\" and, on the very next line, \":\" -> \".\" .\n//    These two lines are separated by a <w:br/> (line break), which shows\n//    up as \"\\u000b\" in the Word.js text model, so searching across it lets\n//    us target this exact \"PA\" + \":\" pair without touching the other lone\n//    \":\" that appears later in the same paragraph (in the \"| :\" line).\nconst paColonResults = body.search(\"PA\\u000b:\", { matchCase: true });\npaColonResults.load(\"items\");\nawait context.sync();\npaColonResults.items.forEach((r) => {\n  r.insertText(\"fA BT\\u000b.\", Word.InsertLocation.replace);\n});\nawait context.sync();\n\n// 5) \"; |\" -> \";\" (drop the trailing \" |\").\nconst semiResults = body.search(\"; |\", { matchCase: true });\nsemiResults.load(\"items\");\nawait context.sync();\nsemiResults.items.forEach((r) => {\n  r.insertText(\";\", Word.InsertLocation.replace);\n});\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\n# 1) Remove the paragraph that holds the inline picture (Picture 1 / 34.jpg).\nif ($d.InlineShapes.Count -gt 0) {\n    $shape = $d.InlineShapes.Item(1)\n    $picPara = $shape.Range.Paragraphs.Item(1)\n    $picPara.Range.Delete()\n}\n\n# 2) Fix the typo \"contalner\" -> \"container\" in the sentence about the water.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"contalner was slowly heated up.\", $false, $false, $false, $false, $false, $true, 1, $false, \"container was slowly heated up.\", 2)\n\n# 3) \"PA\" -> \"fA BT\" and, on the next line (separated by a manual line break,\n#    \"^l\"), \":\" -> \".\" . Matching across the line break targets this exact\n#    pair and avoids the other standalone \":\" later in the same paragraph\n#    (the \"| :\" line), which must stay untouched.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"PA^l:\", $false, $false, $false, $false, $false, $true, 1, $false, \"fA BT^l.\", 2)\n\n# 4) \"; |\" -> \";\" (drop the trailing \" |\").\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"; |\", $false, $false, $false, $false, $false, $true, 1, $false, \";\", 2)\n"}
